$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '68.856.73'
$ws.Range("E2").Value = '  +1.25%  '
# Row 3
$ws.Range("D3").Value = '3.734.55'
$ws.Range("E3").Value = '  -2.56%  '
# Row 4
$ws.Range("E4").Value = '  +0.10%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '601.04'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.30%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.24'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.36%  '
# Row 7
$ws.Range("D7").Value = '3.733.20'
$ws.Range("E7").Value = '  -2.57%  '
# Row 8
$ws.Range("E8").Value = '  -0.01%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.533'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.03%  '
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.163'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.25%  '
# Row 11
$ws.Range("E11").Value = '  +2.64%  '
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.458'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.71%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '37.94'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.76%  '
# Row 14
$ws.Range("E14").Value = '  -0.91%  '
# Row 15
$ws.Range("D15").Value = '4.360.27'
$ws.Range("E15").Value = '  -2.36%  '
# Row 16
$ws.Range("D16").Value = '3.748.10'
$ws.Range("E16").Value = '  -1.94%  '
# Row 17
$ws.Range("D17").Value = '68.831.00'
$ws.Range("E17").Value = '  +1.27%  '
# Row 18
$ws.Range("E18").Value = '  +0.43%  '
# Row 19
$ws.Range("E19").Value = '  +0.28%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.38'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.61%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '497.12'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.94%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.03'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +9.56%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.723'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.95%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.17'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.30%  '
# Row 25
$ws.Range("E25").Value = '  -3.28%  '
# Row 26
$ws.Range("E26").Value = '  -7.73%  '
# Row 27
$ws.Range("E27").Value = '  +0.26%  '
# Row 28
$ws.Range("E28").Value = '  -1.81%  '
# Row 29
$ws.Range("E29").Value = '  -0.13%  '
# Row 30
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.94'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.64%  '
# Row 31
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.46'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.40%  '
# Row 32
$ws.Range("E32").Value = '  +2.41%  '
# Row 33
$ws.Range("E33").Value = '  -4.20%  '
# Row 34
$ws.Range("D34").Value = '3.881.41'
$ws.Range("E34").Value = '  -2.13%  '
# Row 35
$ws.Range("E35").Value = '  -0.58%  '
# Row 36
$ws.Range("D36").Value = '3.657.80'
$ws.Range("E36").Value = '  -2.84%  '
# Row 37
$ws.Range("E37").Value = '  +0.12%  '
# Row 38
$ws.Range("E38").Value = '  -0.86%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.84'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.59%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.132'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.25%  '
# Row 41
$ws.Range("E41").Value = '  -1.81%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '434.00'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.31%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '48.86'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.89%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.97'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.20%  '
# Row 45
$ws.Range("E45").Value = '  -1.92%  '
# Row 46
$ws.Range("E46").Value = '  +0.29%  '
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '40.62'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.86%  '
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '141.36'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.20%  '
# Row 50
$ws.Range("E50").Value = '  +0.29%  '
# Row 51
$ws.Range("D51").Value = '2.741.57'
$ws.Range("E51").Value = '  -3.78%  '
